$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQ_Metrics")

# Remove the "Mixed Cases" metric row (abbreviation mxCase / label "Mixed Cases" / value 18)
$ws.Rows.Item(22).Delete()

# Execution Time value changed from 0.08 to 0.06 (row shifted up after the delete above,
# exe_time now lives on row 27 instead of 28)
$ws.Range("C27").Value = "0.06"

# Clear the "Version of Used Diagnosis List" value (dx_list), keep the row but blank its Value cell
$ws.Range("C32").Value = ""

# Remove the now-undocumented "Encounter Class" row entirely (enctr_class)
$ws.Rows.Item(33).Delete()
